$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Top Gainers")

# Row 2
$ws.Cells.Item(2, 3).Value = 10.7747
$ws.Cells.Item(2, 4).Value = 19.5493
$ws.Cells.Item(2, 5).Value = 26.6812

# Row 3
$ws.Cells.Item(3, 3).Value = 10.5584
$ws.Cells.Item(3, 4).Value = 15.8511
$ws.Cells.Item(3, 5).Value = 22.7734

# Row 4
$ws.Cells.Item(4, 3).Value = 10.2398
$ws.Cells.Item(4, 4).Value = 7.2922
$ws.Cells.Item(4, 5).Value = -7.5409

# Row 5
$ws.Cells.Item(5, 3).Value = 10.2142
$ws.Cells.Item(5, 4).Value = 10.4202
$ws.Cells.Item(5, 5).Value = 24.4492

# Row 6
$ws.Cells.Item(6, 3).Value = 10.1828
$ws.Cells.Item(6, 4).Value = 17.9503
$ws.Cells.Item(6, 5).Value = 31.3877

# Row 7
$ws.Cells.Item(7, 3).Value = 9.5839
$ws.Cells.Item(7, 4).Value = 11.9004
$ws.Cells.Item(7, 5).Value = 28.3762

# Row 8
$ws.Cells.Item(8, 3).Value = 7.2187
$ws.Cells.Item(8, 4).Value = 12.4037
$ws.Cells.Item(8, 5).Value = 14.5619

# Row 9
$ws.Cells.Item(9, 3).Value = 6.7955
$ws.Cells.Item(9, 4).Value = 3.3097
$ws.Cells.Item(9, 5).Value = 15.9254

# Row 10
$ws.Cells.Item(10, 3).Value = 6.2282
$ws.Cells.Item(10, 4).Value = 10.5073
$ws.Cells.Item(10, 5).Value = 11.6033

# Row 11
$ws.Cells.Item(11, 3).Value = 6.2186
$ws.Cells.Item(11, 4).Value = 5.2006
$ws.Cells.Item(11, 5).Value = -23.036

# Row 12
$ws.Cells.Item(12, 3).Value = 5.9154
$ws.Cells.Item(12, 4).Value = 6.2004
$ws.Cells.Item(12, 5).Value = 7.9427

# Row 13
$ws.Cells.Item(13, 2).Value = "PDSL"
$ws.Cells.Item(13, 3).Value = 5.8247
$ws.Cells.Item(13, 4).Value = 9.1775
$ws.Cells.Item(13, 5).Value = 15.3544

# Row 14
$ws.Cells.Item(14, 2).Value = "NETWEB"
$ws.Cells.Item(14, 3).Value = 5.6394
$ws.Cells.Item(14, 4).Value = 11.5251
$ws.Cells.Item(14, 5).Value = 13.622

# Row 15
$ws.Cells.Item(15, 2).Value = "VENKEYS"
$ws.Cells.Item(15, 3).Value = 5.6348
$ws.Cells.Item(15, 4).Value = 6.2654
$ws.Cells.Item(15, 5).Value = 3.9541

# Row 16
$ws.Cells.Item(16, 2).Value = "BLSE"
$ws.Cells.Item(16, 3).Value = 5.5497
$ws.Cells.Item(16, 4).Value = 4.5146
$ws.Cells.Item(16, 5).Value = -1.6867

# Row 17
$ws.Cells.Item(17, 3).Value = 5.4891
$ws.Cells.Item(17, 4).Value = 7.8484
$ws.Cells.Item(17, 5).Value = 6.8155

# Row 18
$ws.Cells.Item(18, 3).Value = 5.0438
$ws.Cells.Item(18, 4).Value = 3.6555
$ws.Cells.Item(18, 5).Value = 8.963900000000001

# Row 22
$ws.Cells.Item(22, 2).Value = "VIMTALABS"
$ws.Cells.Item(22, 3).Value = 4.9543
$ws.Cells.Item(22, 4).Value = 5.1339
$ws.Cells.Item(22, 5).Value = 0

# Row 23
$ws.Cells.Item(23, 2).Value = "EUROPRATIK"
$ws.Cells.Item(23, 3).Value = 4.8382
$ws.Cells.Item(23, 4).Value = 10.7403
$ws.Cells.Item(23, 5).Value = 27.7807

# Row 24
$ws.Cells.Item(24, 3).Value = 4.6375
$ws.Cells.Item(24, 4).Value = 11.1847
$ws.Cells.Item(24, 5).Value = 7.6829

# Row 26
$ws.Cells.Item(26, 3).Value = 4.4896
$ws.Cells.Item(26, 4).Value = 5.7115
$ws.Cells.Item(26, 5).Value = 26.1316

# Row 28
$ws.Cells.Item(28, 3).Value = 4.1606
$ws.Cells.Item(28, 4).Value = 8.3079
$ws.Cells.Item(28, 5).Value = 32.3809

# Row 29
$ws.Cells.Item(29, 3).Value = 4.1281
$ws.Cells.Item(29, 4).Value = 4.6754
$ws.Cells.Item(29, 5).Value = -1.6056

# Row 31
$ws.Cells.Item(31, 2).Value = "ALICON"
$ws.Cells.Item(31, 3).Value = 3.9434
$ws.Cells.Item(31, 4).Value = 10.1053
$ws.Cells.Item(31, 5).Value = 15.5745

# Row 32
$ws.Cells.Item(32, 2).Value = "AHLUCONT"
$ws.Cells.Item(32, 3).Value = 3.8297
$ws.Cells.Item(32, 4).Value = 2.6282
$ws.Cells.Item(32, 5).Value = -4.6523

# Row 33
$ws.Cells.Item(33, 2).Value = "DBL"
$ws.Cells.Item(33, 3).Value = 3.8019
$ws.Cells.Item(33, 4).Value = 4.8995
$ws.Cells.Item(33, 5).Value = 5.9869

# Row 34
$ws.Cells.Item(34, 2).Value = "CANBK"
$ws.Cells.Item(34, 3).Value = 3.7897
$ws.Cells.Item(34, 4).Value = 6.3246
$ws.Cells.Item(34, 5).Value = 8.0349

# Row 35
$ws.Cells.Item(35, 2).Value = "SHANTIGOLD"
$ws.Cells.Item(35, 3).Value = 3.7535
$ws.Cells.Item(35, 4).Value = 11.056
$ws.Cells.Item(35, 5).Value = 3.6467

# Row 36
$ws.Cells.Item(36, 2).Value = "SAPPHIRE"
$ws.Cells.Item(36, 3).Value = 3.7399
$ws.Cells.Item(36, 4).Value = 5.5691
$ws.Cells.Item(36, 5).Value = 2.9101

# Row 37
$ws.Cells.Item(37, 2).Value = "BLISSGVS"
$ws.Cells.Item(37, 3).Value = 3.6961
$ws.Cells.Item(37, 4).Value = 3.0456
$ws.Cells.Item(37, 5).Value = 3.4173

# Row 38
$ws.Cells.Item(38, 2).Value = "ABREL"
$ws.Cells.Item(38, 3).Value = 3.6913
$ws.Cells.Item(38, 4).Value = 12.0364
$ws.Cells.Item(38, 5).Value = 11.5722

# Row 39
$ws.Cells.Item(39, 2).Value = "BLUEDART"
$ws.Cells.Item(39, 3).Value = 3.6823
$ws.Cells.Item(39, 4).Value = 22.4328
$ws.Cells.Item(39, 5).Value = 19.481

# Row 40
$ws.Cells.Item(40, 2).Value = "SKYGOLD"
$ws.Cells.Item(40, 3).Value = 3.6475
$ws.Cells.Item(40, 4).Value = -0.9026
$ws.Cells.Item(40, 5).Value = 37.6258

# Row 41
$ws.Cells.Item(41, 2).Value = "JKTYRE"
$ws.Cells.Item(41, 3).Value = 3.619
$ws.Cells.Item(41, 4).Value = 6.6877
$ws.Cells.Item(41, 5).Value = 22.8855

# Row 42
$ws.Cells.Item(42, 2).Value = "REDTAPE"
$ws.Cells.Item(42, 3).Value = 3.6163
$ws.Cells.Item(42, 4).Value = 3.5157
$ws.Cells.Item(42, 5).Value = -3.4127

# Row 43
$ws.Cells.Item(43, 2).Value = "PFOCUS"
$ws.Cells.Item(43, 3).Value = 3.5673
$ws.Cells.Item(43, 4).Value = 0.846
$ws.Cells.Item(43, 5).Value = 2.3076

# Row 44
$ws.Cells.Item(44, 3).Value = 3.5595
$ws.Cells.Item(44, 4).Value = 6.8243
$ws.Cells.Item(44, 5).Value = -0.4122

# Row 45
$ws.Cells.Item(45, 3).Value = 3.4774
$ws.Cells.Item(45, 4).Value = 3.6526
$ws.Cells.Item(45, 5).Value = 19.4092

# Row 46
$ws.Cells.Item(46, 2).Value = "GMMPFAUDLR"
$ws.Cells.Item(46, 3).Value = 3.4399
$ws.Cells.Item(46, 4).Value = 7.9228
$ws.Cells.Item(46, 5).Value = 20.3171

# Row 47
$ws.Cells.Item(47, 2).Value = "VSTIND"
$ws.Cells.Item(47, 3).Value = 3.3069
$ws.Cells.Item(47, 4).Value = 3.8088
$ws.Cells.Item(47, 5).Value = 3.3469

# Row 48
$ws.Cells.Item(48, 2).Value = "BGRENERGY"
$ws.Cells.Item(48, 3).Value = 3.2153
$ws.Cells.Item(48, 4).Value = -6.0917
$ws.Cells.Item(48, 5).Value = 74.82810000000001

# Row 49
$ws.Cells.Item(49, 3).Value = 3.1611
$ws.Cells.Item(49, 4).Value = -1.2121
$ws.Cells.Item(49, 5).Value = 8.818899999999999

# Row 51
$ws.Cells.Item(51, 2).Value = "FIVESTAR"
$ws.Cells.Item(51, 3).Value = 3.132
$ws.Cells.Item(51, 4).Value = 16.2011
$ws.Cells.Item(51, 5).Value = 16.2877

# Row 52
$ws.Cells.Item(52, 2).Value = "RSYSTEMS"
$ws.Cells.Item(52, 3).Value = 3.0285
$ws.Cells.Item(52, 4).Value = 4.185
$ws.Cells.Item(52, 5).Value = 6.5245

# Row 53
$ws.Cells.Item(53, 2).Value = "ASHOKA"
$ws.Cells.Item(53, 3).Value = 3.0178
$ws.Cells.Item(53, 4).Value = 4.5409
$ws.Cells.Item(53, 5).Value = 7.1722

# Row 54
$ws.Cells.Item(54, 2).Value = "BLS"
$ws.Cells.Item(54, 3).Value = 3.0112
$ws.Cells.Item(54, 4).Value = -0.0314
$ws.Cells.Item(54, 5).Value = -1.2876

# Row 55
$ws.Cells.Item(55, 2).Value = "SUNDROP"
$ws.Cells.Item(55, 3).Value = 2.9786
$ws.Cells.Item(55, 4).Value = 2.8713
$ws.Cells.Item(55, 5).Value = 0.9714

# Row 56
$ws.Cells.Item(56, 2).Value = "MRPL"
$ws.Cells.Item(56, 3).Value = 2.9106
$ws.Cells.Item(56, 4).Value = 12.9034
$ws.Cells.Item(56, 5).Value = 23.5485

# Row 57
$ws.Cells.Item(57, 2).Value = "PSPPROJECT"
$ws.Cells.Item(57, 3).Value = 2.9047
$ws.Cells.Item(57, 4).Value = 16.9828
$ws.Cells.Item(57, 5).Value = 23.4144

# Row 58
$ws.Cells.Item(58, 2).Value = "BPCL"
$ws.Cells.Item(58, 3).Value = 2.8727
$ws.Cells.Item(58, 4).Value = 8.3674
$ws.Cells.Item(58, 5).Value = 5.4321

# Row 59
$ws.Cells.Item(59, 2).Value = "VOLTAMP"
$ws.Cells.Item(59, 3).Value = 2.8508
$ws.Cells.Item(59, 4).Value = 2.7297
$ws.Cells.Item(59, 5).Value = 2.467

# Row 60
$ws.Cells.Item(60, 2).Value = "CENTRUM"
$ws.Cells.Item(60, 3).Value = 2.7576
$ws.Cells.Item(60, 4).Value = 1.3752
$ws.Cells.Item(60, 5).Value = 0.7128

# Row 62
$ws.Cells.Item(62, 2).Value = "FEDFINA"
$ws.Cells.Item(62, 3).Value = 2.6868
$ws.Cells.Item(62, 4).Value = 3.7907
$ws.Cells.Item(62, 5).Value = -4.9605

# Row 63
$ws.Cells.Item(63, 2).Value = "CIFL"
$ws.Cells.Item(63, 3).Value = 2.6461
$ws.Cells.Item(63, 4).Value = 2.2003
$ws.Cells.Item(63, 5).Value = 2.1412

# Row 64
$ws.Cells.Item(64, 2).Value = "GANESHCP"
$ws.Cells.Item(64, 3).Value = 2.6027
$ws.Cells.Item(64, 4).Value = 2.0786
$ws.Cells.Item(64, 5).Value = 1.631

# Row 65
$ws.Cells.Item(65, 2).Value = "OIL"
$ws.Cells.Item(65, 3).Value = 2.5922
$ws.Cells.Item(65, 4).Value = 2.8367
$ws.Cells.Item(65, 5).Value = 4.2533

# Row 66
$ws.Cells.Item(66, 2).Value = "REFEX"
$ws.Cells.Item(66, 3).Value = 2.5247
$ws.Cells.Item(66, 4).Value = -0.0413
$ws.Cells.Item(66, 5).Value = 1.9352

# Row 67
$ws.Cells.Item(67, 3).Value = 2.5239
$ws.Cells.Item(67, 4).Value = 9.4733
$ws.Cells.Item(67, 5).Value = 18.6546

# Row 68
$ws.Cells.Item(68, 2).Value = "MFSL"
$ws.Cells.Item(68, 3).Value = 2.5217
$ws.Cells.Item(68, 4).Value = 2.5758
$ws.Cells.Item(68, 5).Value = -1.193

# Row 69
$ws.Cells.Item(69, 2).Value = "CENTUM"
$ws.Cells.Item(69, 3).Value = 2.5191
$ws.Cells.Item(69, 4).Value = 3.375
$ws.Cells.Item(69, 5).Value = -2.1057

# Row 70
$ws.Cells.Item(70, 2).Value = "INOXGREEN"
$ws.Cells.Item(70, 3).Value = 2.4718
$ws.Cells.Item(70, 4).Value = 10.4611
$ws.Cells.Item(70, 5).Value = 33.7825

# Row 71
$ws.Cells.Item(71, 2).Value = "OBEROIRLTY"
$ws.Cells.Item(71, 3).Value = 2.4688
$ws.Cells.Item(71, 4).Value = 3.2825
$ws.Cells.Item(71, 5).Value = 10.9728

# Row 72
$ws.Cells.Item(72, 2).Value = "INDORAMA"
$ws.Cells.Item(72, 3).Value = 2.4344
$ws.Cells.Item(72, 4).Value = 4.8963
$ws.Cells.Item(72, 5).Value = 16.2641

# Row 73
$ws.Cells.Item(73, 3).Value = 2.4251
$ws.Cells.Item(73, 4).Value = 1.7731
$ws.Cells.Item(73, 5).Value = 5.0657

# Row 74
$ws.Cells.Item(74, 2).Value = "JKLAKSHMI"
$ws.Cells.Item(74, 3).Value = 2.379
$ws.Cells.Item(74, 4).Value = 4.364
$ws.Cells.Item(74, 5).Value = 1.3879

# Row 75
$ws.Cells.Item(75, 2).Value = "SPANDANA"
$ws.Cells.Item(75, 3).Value = 2.3719
$ws.Cells.Item(75, 4).Value = 3.8625
$ws.Cells.Item(75, 5).Value = 2.4822

# Row 76
$ws.Cells.Item(76, 2).Value = "SDBL"
$ws.Cells.Item(76, 3).Value = 2.313
$ws.Cells.Item(76, 4).Value = 0.8617
$ws.Cells.Item(76, 5).Value = 6.4381
